$wb = $excel.ActiveWorkbook

# 1. Insert the new worksheet "dadosPessoaisCompra" right after "dadosDeCadastro"
#    (i.e. right before "dadosDeAcesso").
$wsCadastro = $wb.Worksheets.Item("dadosDeCadastro")
$wsCompra = $wb.Worksheets.Add($null, $wsCadastro)
$wsCompra.Name = "dadosPessoaisCompra"

# 2. Fill in column A (id header reuses the existing "id" string; the
#    ID_0012 value is entered first, matching the authored edit order).
$wsCompra.Range("A1").Value = "id"
$wsCompra.Range("A2").Value = "ID_0012"

# 3. Fill in the header row.
$wsCompra.Range("B1").Value = "nome"
$wsCompra.Range("C1").Value = "país"
$wsCompra.Range("D1").Value = "cidade"
$wsCompra.Range("E1").Value = "cartão"
$wsCompra.Range("F1").Value = "mês"
$wsCompra.Range("G1").Value = "ano"

# 4. Fill in the rest of the data row.
$wsCompra.Range("B2").Value = "André"
$wsCompra.Range("C2").Value = "brasil"
$wsCompra.Range("D2").Value = "santo andré"

$wsCompra.Range("E2:G2").NumberFormat = "@"
$wsCompra.Range("E2").Value = '"1234567898765"'
$wsCompra.Range("F2").Value = '"8"'
$wsCompra.Range("G2").Value = '"1998"'
$wsCompra.Range("G2").Font.Underline = $true

# 5. Column widths (best-fit-ish) for columns D and E.
$wsCompra.Columns.Item(4).ColumnWidth = 10.6
$wsCompra.Columns.Item(5).ColumnWidth = 18.2

# 6. Page setup to match the other sheets.
$wsCompra.PageSetup.PaperSize = 9
$wsCompra.PageSetup.Orientation = 1

# 7. Select G3 (just past the filled data) and make this the active sheet/tab.
$wsCompra.Activate()
$wsCompra.Range("G3").Select()

# 8. Add a new row (ID_0012) to the "dadosDeAcesso" sheet, following the existing pattern.
$wsAcesso = $wb.Worksheets.Item("dadosDeAcesso")
$wsAcesso.Cells.Item(10, 1).Value = "ID_0012"
$wsAcesso.Cells.Item(10, 2).Value = "André Automatizador"
$wsAcesso.Cells.Item(10, 3).Value = "sem email"
$wsAcesso.Cells.Item(10, 4).Value = "automacaoteste"

# 9. Re-activate dadosPessoaisCompra as the final active tab (Cells.Item above may have
#    shifted the active sheet), and keep dadosDeAcesso's own selection on A10.
$wsAcesso.Range("A10").Select()
$wsCompra.Activate()
$wsCompra.Range("G3").Select()

Write-Output "done"
